$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97, pushing the existing rows 97-107 down to 98-108.
$ws.Rows.Item(97).Insert()

$row = 97
$ws.Cells.Item($row, 1).Value2 = 10
$ws.Cells.Item($row, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value2 = "La Araucanía"
$ws.Cells.Item($row, 4).Value2 = 45223
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value2 = 9
$ws.Cells.Item($row, 6).Value2 = 100112042
$ws.Cells.Item($row, 7).Value2 = "Locoto"
$ws.Cells.Item($row, 8).Value2 = "Sin especificar"
$ws.Cells.Item($row, 9).Value2 = "Primera"
$ws.Cells.Item($row, 10).Value2 = 120
$ws.Cells.Item($row, 11).Value2 = 3800
$ws.Cells.Item($row, 12).Value2 = 3800
$ws.Cells.Item($row, 13).Value2 = 3800
$ws.Cells.Item($row, 14).Value2 = "$/kilo"
$ws.Cells.Item($row, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value2 = 3800
$ws.Cells.Item($row, 17).Value2 = 1
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"
